$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Ohio State', 28.1),
    @('Indiana', 24.2),
    @('Notre Dame', 24.1),
    @('Alabama', 23.7),
    @('Texas', 22.5),
    @('Oregon', 22.1),
    @('Georgia', 21.4),
    @('Texas Tech', 21.1),
    @('Miami', 21.1),
    @('Texas A&M', 19.5),
    @('USC', 18.8),
    @('Utah', 18.3),
    @('LSU', 17.5),
    @('Ole Miss', 16.6),
    @('Oklahoma', 16.4),
    @('Michigan', 15.6),
    @('Tennessee', 15.5),
    @('Missouri', 15),
    @('Vanderbilt', 14.9),
    @('Penn State', 14.8),
    @('Iowa', 13.9),
    @('Clemson', 13.5),
    @('Florida', 13.4),
    @('Washington', 13.4),
    @('BYU', 13.2),
    @('Auburn', 13.1),
    @('Nebraska', 12.8),
    @('Florida State', 12.1),
    @('Illinois', 11.4),
    @('Louisville', 10.2),
    @('South Carolina', 10.1),
    @('Iowa State', 10),
    @('TCU', 9.9),
    @('Georgia Tech', 9.7),
    @('Duke', 9.4),
    @('Cincinnati', 9.1),
    @('Pittsburgh', 8.9),
    @('Arkansas', 8.6),
    @('SMU', 8.4),
    @('South Florida', 8.4),
    @('Kansas', 8.3),
    @('Virginia', 8.2),
    @('Kansas State', 8.2),
    @('Memphis', 7.7),
    @('Arizona State', 7.6),
    @('Baylor', 7.6),
    @('Mississippi State', 6.6),
    @('Boise State', 6.2),
    @('Arizona', 5.3),
    @('Rutgers', 4.6),
    @('Maryland', 4.2),
    @('Houston', 4.1),
    @('Kentucky', 4),
    @('Colorado', 4),
    @('NC State', 3.7),
    @('UCF', 3.4),
    @('Tulane', 3.2),
    @('Northwestern', 3),
    @('East Carolina', 2.4),
    @('Minnesota', 2.4),
    @('Old Dominion', 2.3),
    @('UCLA', 1.8),
    @('Toledo', 1.6),
    @('James Madison', 1.5),
    @('Wisconsin', 0.5),
    @('Syracuse', 0.4),
    @('North Texas', 0.1),
    @('Wake Forest', 0),
    @('UNLV', -0.7),
    @('San Diego State', -0.7),
    @('UTSA', -0.9),
    @('Virginia Tech', -1),
    @('Michigan State', -1),
    @('Navy', -2.2),
    @('West Virginia', -2.7),
    @('UConn', -3),
    @('California', -3.2),
    @('Texas State', -3.7),
    @('Ohio', -3.9),
    @('Purdue', -4.1),
    @('Army', -5.2),
    @('Boston College', -5.5),
    @('Troy', -5.7),
    @('Marshall', -5.9),
    @('Miami (OH)', -6.2),
    @('Stanford', -6.4),
    @('North Carolina', -6.6),
    @('Washington State', -6.8),
    @('New Mexico', -6.9),
    @('Fresno State', -7),
    @('Western Michigan', -7.5),
    @('Louisiana Tech', -7.6),
    @('Bowling Green', -7.7),
    @('Utah State', -7.8),
    @('Hawai''i', -7.9),
    @('Air Force', -7.9),
    @('Temple', -8),
    @('Colorado State', -8.1),
    @('Southern Miss', -8.5),
    @('Western Kentucky', -8.5),
    @('San José State', -8.5),
    @('Liberty', -9.4),
    @('Wyoming', -9.6),
    @('South Alabama', -9.6),
    @('App State', -9.8),
    @('Oregon State', -9.9),
    @('Louisiana', -10.4),
    @('Buffalo', -10.6),
    @('Jacksonville State', -11.5),
    @('Georgia Southern', -11.6),
    @('Kennesaw State', -11.8),
    @('Delaware', -12.2),
    @('Florida Atlantic', -12.5),
    @('Oklahoma State', -12.8),
    @('Arkansas State', -13.3),
    @('Northern Illinois', -13.4),
    @('Coastal Carolina', -14.3),
    @('Central Michigan', -15.8),
    @('Missouri State', -15.8),
    @('Rice', -16),
    @('Eastern Michigan', -16.5),
    @('UL Monroe', -16.6),
    @('UTEP', -16.7),
    @('Tulsa', -17.5),
    @('UAB', -17.6),
    @('New Mexico State', -18.4),
    @('Florida International', -18.7),
    @('Georgia State', -19),
    @('Nevada', -19.1),
    @('Charlotte', -20.6),
    @('Middle Tennessee', -20.7),
    @('Akron', -21),
    @('Sam Houston', -21),
    @('Ball State', -21.6),
    @('Kent State', -22.3),
    @('Massachusetts', -28.5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Activate()
$ws.Range("G124").Select()